# COM-interop edit script: refresh COVID-19 "paises" snapshot
# (matches upstream commit "Update countries & provincias Spain")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "last updated" timestamp shown in the title cell (A1)
$ws.Range("A1").Value = 'Datos actualizados a 26 de Marzo de 2020 a las 21:42'

# 2) Refresh the per-country statistics with the newer snapshot.
#    Rows whose case counts changed rank also get column A (the
#    country name) rewritten so the table stays sorted by
#    "Casos totales" descending, matching the source data refresh.
$rowsData = @{
    4 = @('Estados Unidos'; 81384; 13173; 1864; 78347; 2112; 146; 1173)
    5 = @('China'; 81285; 67; 74051; 3947; 1235; 6; 3287)
    8 = @('Alemania'; 43646; 6323; 5673; 37711; 23; 56; 262)
    11 = @('Suiza'; 11811; 914; 131; 11489; 141; 38; 191)
    21 = @('Brasil'; 2915; 361; 6; 2832; 18; 18; 77)
    22 = @('Suecia'; 2840; 314; 16; 2753; 176; 9; 71)
    23 = @('Australia'; 2806; 130; 170; 2623; 11; 2; 13)
    24 = @('Israel'; 2693; 324; 68; 2617; 46; 3; 8)
    55 = @('Egipto'; 495; 39; 102; 369; 0; 3; 24)
    56 = @('Colombia'; 491; 21; 8; 477; 0; 2; 6)
    57 = @('Republica Dominicana'; 488; 96; 3; 475; 0; 0; 10)
    58 = @('Mexico'; 475; 70; 4; 465; 1; 1; 6)
    59 = @('Barein'; 458; 39; 204; 250; 1; 0; 4)
    60 = @('Serbia'; 457; 73; 15; 435; 21; 3; 7)
    74 = @('Costa Rica'; 231; 30; 2; 227; 5; 0; 2)
    83 = @('Ucrania'; 196; 51; 1; 190; 0; 0; 5)
    84 = @('Bosnia y Herzegovina'; 189; 13; 2; 184; 1; 0; 3)
    85 = @('Moldavia'; 177; 28; 2; 174; 28; 0; 1)
    86 = @('Albania'; 174; 28; 17; 151; 3; 1; 6)
    101 = @('Costa de Marfil'; 96; 16; 3; 93; 0; 0; 0)
    102 = @('Camboya'; 96; 0; 10; 86; 1; 0; 0)
    103 = @('Afganistan'; 94; 10; 2; 90; 0; 0; 2)
    104 = @('Bielorrusia'; 86; 0; 29; 57; 2; 0; 0)
    105 = @('Estado de Palestina'; 84; 13; 17; 66; 0; 0; 1)
    106 = @('Mauricio'; 81; 33; 0; 79; 1; 0; 2)
    121 = @('Bolivia'; 43; 11; 0; 43; 0; 0; 0)
    122 = @('Ruanda'; 41; 0; 0; 41; 0; 0; 0)
    123 = @('Paraguay'; 41; 4; 0; 38; 1; 0; 3)
    142 = @('Nueva Caledonia'; 14; 0; 0; 14; 0; 0; 0)
    143 = @('Uganda'; 14; 0; 0; 14; 0; 0; 0)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}
